$wb = $excel.ActiveWorkbook

# Rename the "Include from Ferlab.bio CodeS" sheet to "Include #0"
$wsInclude = $wb.Worksheets.Item("Include from Ferlab.bio CodeS")
$wsInclude.Name = "Include #0"

# Metadata sheet updates
$wsMeta = $wb.Worksheets.Item("Metadata")

# Date value changed
$wsMeta.Range("B8").Value = "2024-10-02T15:04:17+00:00"

# Contact value changed
$wsMeta.Range("B10").Value = "Ferlab.bio (http://example.org/example-publisher)"

# Insert a new row for "Jurisdiction" before "Description" (which is currently row 11)
$wsMeta.Rows("11").Insert()
$wsMeta.Range("A11").Value = "Jurisdiction"
$wsMeta.Range("B11").Value = ""

# Copy the formatting of the surrounding data rows onto the new row
$wsMeta.Range("A10:B10").Copy()
$wsMeta.Range("A11:B11").PasteSpecial(-4122)
$wsMeta.Range("A11").Value = "Jurisdiction"
$wsMeta.Range("B11").Value = ""
